$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 1938
$ws.Cells.Item(2, 2).Value = "Lavínia Souza"
$ws.Cells.Item(2, 3).Value = "Operacoes"
$ws.Cells.Item(2, 4).Value = "Viagem de negocios"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 45088
$ws.Cells.Item(2, 7).Value = 9842.9

# Row 3
$ws.Cells.Item(3, 1).Value = 48328
$ws.Cells.Item(3, 2).Value = "Nathan Azevedo"
$ws.Cells.Item(3, 3).Value = "Marketing"
$ws.Cells.Item(3, 4).Value = "Consulta medica"
$ws.Cells.Item(3, 5).Value = 6
$ws.Cells.Item(3, 6).Value = 45079
$ws.Cells.Item(3, 7).Value = 4280.77

# Row 4
$ws.Cells.Item(4, 1).Value = 32590
$ws.Cells.Item(4, 2).Value = "Bryan Pimenta"
$ws.Cells.Item(4, 3).Value = "P&D"
$ws.Cells.Item(4, 4).Value = "Viagem de negocios"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 45089
$ws.Cells.Item(4, 7).Value = 7330.7

# Row 5
$ws.Cells.Item(5, 1).Value = 89196
$ws.Cells.Item(5, 2).Value = "Maria Laura Garcia"
$ws.Cells.Item(5, 3).Value = "TI"
$ws.Cells.Item(5, 4).Value = "Problemas pessoais"
$ws.Cells.Item(5, 5).Value = 4
$ws.Cells.Item(5, 6).Value = 45094
$ws.Cells.Item(5, 7).Value = 9084.25

# Row 6
$ws.Cells.Item(6, 1).Value = 31566
$ws.Cells.Item(6, 2).Value = "Sr. João Vitor Guerra"
$ws.Cells.Item(6, 3).Value = "Vendas"
$ws.Cells.Item(6, 4).Value = "Outros"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 45094
$ws.Cells.Item(6, 7).Value = 6626.38

# Row 7
$ws.Cells.Item(7, 1).Value = 56284
$ws.Cells.Item(7, 2).Value = "Isaac Silveira"
$ws.Cells.Item(7, 3).Value = "Juridico"
$ws.Cells.Item(7, 4).Value = "Viagem de negocios"
$ws.Cells.Item(7, 5).Value = 7
$ws.Cells.Item(7, 6).Value = 45088
$ws.Cells.Item(7, 7).Value = 6916.51

# Row 8
$ws.Cells.Item(8, 1).Value = 64342
$ws.Cells.Item(8, 2).Value = "Larissa Azevedo"
$ws.Cells.Item(8, 3).Value = "Operacoes"
$ws.Cells.Item(8, 4).Value = "Doenca"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 45083
$ws.Cells.Item(8, 7).Value = 8233.57

# Row 9
$ws.Cells.Item(9, 1).Value = 63813
$ws.Cells.Item(9, 2).Value = "Cauê Cavalcanti"
$ws.Cells.Item(9, 3).Value = "Juridico"
$ws.Cells.Item(9, 4).Value = "Viagem de negocios"
$ws.Cells.Item(9, 5).Value = 5
$ws.Cells.Item(9, 6).Value = 45105
$ws.Cells.Item(9, 7).Value = 2086.99

# Row 10
$ws.Cells.Item(10, 1).Value = 70543
$ws.Cells.Item(10, 2).Value = "Emanuel Cunha"
$ws.Cells.Item(10, 5).Value = 4
$ws.Cells.Item(10, 6).Value = 45095
$ws.Cells.Item(10, 7).Value = 4666.57

# Row 11
$ws.Cells.Item(11, 1).Value = 69266
$ws.Cells.Item(11, 2).Value = "Cauã Pires"
$ws.Cells.Item(11, 3).Value = "Vendas"
$ws.Cells.Item(11, 4).Value = "Outros"
$ws.Cells.Item(11, 5).Value = 8
$ws.Cells.Item(11, 6).Value = 45097
$ws.Cells.Item(11, 7).Value = 8982.17
